$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - add May (G) and June (H) values
$ws.Range("G2").Value = 50.06
$ws.Range("H2").Value = 47.66

# Row 3 - update April (F) value, add May (G) and June (H) values
$ws.Range("F3").Value = 81.01000000000001
$ws.Range("G3").Value = 80.43000000000001
$ws.Range("H3").Value = 68.48

# Row 4 - add May (G) and June (H) values
$ws.Range("G4").Value = 84.06
$ws.Range("H4").Value = 88.38

# Row 5 - add May (G) and June (H) values
$ws.Range("G5").Value = 66.31999999999999
$ws.Range("H5").Value = 73.01000000000001
